$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '25.669.37'
$ws.Range("E2").Value = '  -3.69%  '
$ws.Range("D3").Value = '1.740.53'
$ws.Range("E3").Value = '  -5.66%  '
$ws.Range("D4").Value = "'1.003"
$ws.Range("E4").Value = '  +0.26%  '
$ws.Range("D5").Value = "'237.17"
$ws.Range("E5").Value = '  -8.84%  '
$ws.Range("D7").Value = "'0.4924"
$ws.Range("E7").Value = '  -6.79%  '
$ws.Range("D8").Value = "'41.63"
$ws.Range("E8").Value = '  -7.42%  '
$ws.Range("D9").Value = "'0.2373"
$ws.Range("E9").Value = '  -24.99%  '
$ws.Range("D10").Value = "'0.05948"
$ws.Range("E10").Value = '  -12.61%  '
$ws.Range("D11").Value = '1.739.81'
$ws.Range("E11").Value = '  -5.73%  '
$ws.Range("D12").Value = "'0.06834"
$ws.Range("E12").Value = '  -12.27%  '
$ws.Range("D13").Value = "'14.54"
$ws.Range("E13").Value = '  -23.37%  '
$ws.Range("D14").Value = "'4.447"
$ws.Range("E14").Value = '  -11.38%  '
$ws.Range("D15").Value = "'76.87"
$ws.Range("E15").Value = '  -12.92%  '
$ws.Range("D16").Value = "'0.5699"
$ws.Range("E16").Value = '  -27.68%  '
$ws.Range("D17").Value = "'1.003"
$ws.Range("E17").Value = '  +0.32%  '
$ws.Range("D18").Value = "'1.002"
$ws.Range("E18").Value = '  +0.11%  '
$ws.Range("D19").Value = '25.721.60'
$ws.Range("E19").Value = '  -3.54%  '
$ws.Range("D20").Value = "'11.42"
$ws.Range("E20").Value = '  -17.91%  '
$ws.Range("D21").Value = "'0.000006405"
$ws.Range("E21").Value = '  -19.31%  '
$ws.Range("D22").Value = '1.959.27'
$ws.Range("E22").Value = '  -6.02%  '
$ws.Range("D23").Value = "'3.950"
$ws.Range("E23").Value = '  -14.37%  '
$ws.Range("D24").Value = "'5.040"
$ws.Range("E24").Value = '  -15.80%  '
$ws.Range("D25").Value = "'7.753"
$ws.Range("E25").Value = '  -17.03%  '
$ws.Range("D26").Value = "'136.76"
$ws.Range("E26").Value = '  -4.35%  '
$ws.Range("D27").Value = "'1.471"
$ws.Range("E27").Value = '  -12.41%  '
$ws.Range("D28").Value = "'1.826"
$ws.Range("E28").Value = '  -17.99%  '
$ws.Range("D29").Value = "'14.50"
$ws.Range("E29").Value = '  -15.08%  '
$ws.Range("D30").Value = "'100.36"
$ws.Range("E30").Value = '  -9.58%  '
$ws.Range("D31").Value = "'3.768"
$ws.Range("E31").Value = '  -10.51%  '
$ws.Range("D32").Value = "'0.08094"
$ws.Range("E32").Value = '  -7.23%  '
$ws.Range("D33").Value = "'3.346"
$ws.Range("E33").Value = '  -18.27%  '
$ws.Range("D34").Value = "'0.04355"
$ws.Range("E34").Value = '  -10.93%  '
$ws.Range("D35").Value = "'1.001"
$ws.Range("E35").Value = '  +0.11%  '
$ws.Range("D36").Value = "'2.683"
$ws.Range("E36").Value = '  -6.10%  '
$ws.Range("D37").Value = "'1.012"
$ws.Range("E37").Value = '  -11.45%  '
$ws.Range("D38").Value = "'0.6036"
$ws.Range("E38").Value = '  -17.68%  '
$ws.Range("D39").Value = "'2.709"
$ws.Range("E39").Value = '  -12.84%  '
$ws.Range("D40").Value = "'2.060"
$ws.Range("E40").Value = '  -10.20%  '
$ws.Range("D41").Value = "'1.001"
$ws.Range("E41").Value = '  +0.10%  '
$ws.Range("D42").Value = "'103.00"
$ws.Range("E42").Value = '  -6.34%  '
$ws.Range("D43").Value = "'0.01477"
$ws.Range("E43").Value = '  -14.67%  '
$ws.Range("D44").Value = "'0.7775"
$ws.Range("E44").Value = '  -13.78%  '
$ws.Range("D45").Value = "'5.126"
$ws.Range("E45").Value = '  -13.97%  '
$ws.Range("D46").Value = "'0.3751"
$ws.Range("E46").Value = '  -22.34%  '
$ws.Range("D47").Value = "'0.05107"
$ws.Range("E47").Value = '  -12.34%  '
$ws.Range("D48").Value = "'5.963"
$ws.Range("E48").Value = '  -22.76%  '
$ws.Range("D49").Value = "'0.1067"
$ws.Range("E49").Value = '  -14.10%  '
$ws.Range("D50").Value = "'30.23"
$ws.Range("E50").Value = '  -13.31%  '
$ws.Range("D51").Value = "'52.50"
$ws.Range("E51").Value = '  -12.45%  '
